# "Generate Report for Handback"
#
# For both the zh-cn and de-de sheets, the two real source rows (rows 2-3)
# move from "Ready for handoff" to "Handed back: in sync with en-US":
#   - Status (col B) text changes
#   - Latest Target File (col E) / Latest Handback File (col F) are filled
#     in with hyperlinks (pointing at the same targets as the existing
#     Source File Name / Latest Handoff File links)
#   - Latest Handback DateTime (col G) gets a real timestamp instead of the
#     "never happened" sentinel
# Row 4 (.localization-config, "Not to be localized") is unaffected.

$wb = $excel.ActiveWorkbook

# zh-cn sheet
$zhWs = $wb.Worksheets.Item("zh-cn")

$zhWs.Range("B2").Value = "Handed back: in sync with en-US"
$zhWs.Range("E2").Value = "09766ec6-d738-4ddc-8009-71feb9ccd36a.md"
$zhWs.Hyperlinks.Add($zhWs.Range("E2"), "https://github.com/OpenLocalizationTest/oltest/blob/26c520c1310f6696332c13693df59c5373ab030d/e2e/09766ec6-d738-4ddc-8009-71feb9ccd36a.md", "", "", "09766ec6-d738-4ddc-8009-71feb9ccd36a.md")
$zhWs.Range("E2").Font.Underline = 2
$zhWs.Range("E2").Font.Color = 15570276

$zhWs.Range("F2").Value = "09766ec6-d738-4ddc-8009-71feb9ccd36a.4ed35deb9582fba467ef42e739fae176392ff72d.zh-cn.xlf"
$zhWs.Hyperlinks.Add($zhWs.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/96b8ea7e7708cefd9d6d6eaba92615b85141bf18/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/09766ec6-d738-4ddc-8009-71feb9ccd36a.4ed35deb9582fba467ef42e739fae176392ff72d.zh-cn.xlf", "", "", "09766ec6-d738-4ddc-8009-71feb9ccd36a.4ed35deb9582fba467ef42e739fae176392ff72d.zh-cn.xlf")
$zhWs.Range("F2").Font.Underline = 2
$zhWs.Range("F2").Font.Color = 15570276

$zhWs.Range("G2").Value = "2016-02-25 06:19:24"

$zhWs.Range("B3").Value = "Handed back: in sync with en-US"
$zhWs.Range("E3").Value = "73eddf8e-57b7-4f20-af60-4cbd50e8b1f4.md"
$zhWs.Hyperlinks.Add($zhWs.Range("E3"), "https://github.com/OpenLocalizationTest/oltest/blob/26c520c1310f6696332c13693df59c5373ab030d/e2e/73eddf8e-57b7-4f20-af60-4cbd50e8b1f4.md", "", "", "73eddf8e-57b7-4f20-af60-4cbd50e8b1f4.md")
$zhWs.Range("E3").Font.Underline = 2
$zhWs.Range("E3").Font.Color = 15570276

$zhWs.Range("F3").Value = "73eddf8e-57b7-4f20-af60-4cbd50e8b1f4.7dd098297e18f63199724706424b9beead9af350.zh-cn.xlf"
$zhWs.Hyperlinks.Add($zhWs.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/96b8ea7e7708cefd9d6d6eaba92615b85141bf18/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/73eddf8e-57b7-4f20-af60-4cbd50e8b1f4.7dd098297e18f63199724706424b9beead9af350.zh-cn.xlf", "", "", "73eddf8e-57b7-4f20-af60-4cbd50e8b1f4.7dd098297e18f63199724706424b9beead9af350.zh-cn.xlf")
$zhWs.Range("F3").Font.Underline = 2
$zhWs.Range("F3").Font.Color = 15570276

$zhWs.Range("G3").Value = "2016-02-25 06:19:24"

# de-de sheet
$deWs = $wb.Worksheets.Item("de-de")

$deWs.Range("B2").Value = "Handed back: in sync with en-US"
$deWs.Range("E2").Value = "09766ec6-d738-4ddc-8009-71feb9ccd36a.md"
$deWs.Hyperlinks.Add($deWs.Range("E2"), "https://github.com/OpenLocalizationTest/oltest/blob/26c520c1310f6696332c13693df59c5373ab030d/e2e/09766ec6-d738-4ddc-8009-71feb9ccd36a.md", "", "", "09766ec6-d738-4ddc-8009-71feb9ccd36a.md")
$deWs.Range("E2").Font.Underline = 2
$deWs.Range("E2").Font.Color = 15570276

$deWs.Range("F2").Value = "09766ec6-d738-4ddc-8009-71feb9ccd36a.4ed35deb9582fba467ef42e739fae176392ff72d.de-de.xlf"
$deWs.Hyperlinks.Add($deWs.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9385c5c15284bbf465842bd85fb12331b9504f9f/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/09766ec6-d738-4ddc-8009-71feb9ccd36a.4ed35deb9582fba467ef42e739fae176392ff72d.de-de.xlf", "", "", "09766ec6-d738-4ddc-8009-71feb9ccd36a.4ed35deb9582fba467ef42e739fae176392ff72d.de-de.xlf")
$deWs.Range("F2").Font.Underline = 2
$deWs.Range("F2").Font.Color = 15570276

$deWs.Range("G2").Value = "2016-02-25 06:19:46"

$deWs.Range("B3").Value = "Handed back: in sync with en-US"
$deWs.Range("E3").Value = "73eddf8e-57b7-4f20-af60-4cbd50e8b1f4.md"
$deWs.Hyperlinks.Add($deWs.Range("E3"), "https://github.com/OpenLocalizationTest/oltest/blob/26c520c1310f6696332c13693df59c5373ab030d/e2e/73eddf8e-57b7-4f20-af60-4cbd50e8b1f4.md", "", "", "73eddf8e-57b7-4f20-af60-4cbd50e8b1f4.md")
$deWs.Range("E3").Font.Underline = 2
$deWs.Range("E3").Font.Color = 15570276

$deWs.Range("F3").Value = "73eddf8e-57b7-4f20-af60-4cbd50e8b1f4.7dd098297e18f63199724706424b9beead9af350.de-de.xlf"
$deWs.Hyperlinks.Add($deWs.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9385c5c15284bbf465842bd85fb12331b9504f9f/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/73eddf8e-57b7-4f20-af60-4cbd50e8b1f4.7dd098297e18f63199724706424b9beead9af350.de-de.xlf", "", "", "73eddf8e-57b7-4f20-af60-4cbd50e8b1f4.7dd098297e18f63199724706424b9beead9af350.de-de.xlf")
$deWs.Range("F3").Font.Underline = 2
$deWs.Range("F3").Font.Color = 15570276

$deWs.Range("G3").Value = "2016-02-25 06:19:46"

Write-Host "Handback report generated for zh-cn and de-de sheets."
